$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) "Jewel_GoldLoan_Renewal_Transfer" sheet — update the two amount
#    cells (K2, N2) as part of reworking the Jewel Loan Renewal Transfer.
# -----------------------------------------------------------------------
$wsTransfer = $wb.Worksheets.Item("Jewel_GoldLoan_Renewal_Transfer")
$wsTransfer.Range("K2").Value = 195458
$wsTransfer.Range("N2").Value = 329375

# -----------------------------------------------------------------------
# 2) Add the new "GL_Acc_Trans_Cash" test-scenario sheet at the end of
#    the workbook. It mirrors the layout of the existing
#    "Jewel_Loan_Transcharge_Cash" sheet, so clone that sheet and adjust
#    the scenario-name cell.
# -----------------------------------------------------------------------
$wsSrc  = $wb.Worksheets.Item("Jewel_Loan_Transcharge_Cash")
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSrc.Copy($null, $wsLast) | Out-Null

$wsNew = $wb.ActiveSheet
$wsNew.Name = "GL_Acc_Trans_Cash"
$wsNew.Range("A2").Value = "GL_Acc_Transaction_Cash"

# The cloned row kept the source's explicit row height; auto-fit it back
# to the sheet default so no custom height is stored.
$wsNew.Rows.Item(2).AutoFit() | Out-Null

# Selection left on the new sheet before focus moves elsewhere.
$wsNew.Range("D8").Select() | Out-Null

# -----------------------------------------------------------------------
# 3) Leave "Jewel_GoldLoan_Renewal_Transfer" as the active/selected tab,
#    with its selection on N5 (this also clears tabSelected from
#    whichever sheet — "Jewel_Loan_Transcharge_Transfer" / the new sheet
#    — previously held it).
# -----------------------------------------------------------------------
$wsTransfer.Activate()
$wsTransfer.Range("N5").Select() | Out-Null
